# Update Test Data Path
#
# Refresh the sample CA/BM test-data names on the FacilityRegister sheet
# (DIYAS -> KALEB, YOMMIE -> CHAIRUL) and leave that sheet active/selected
# instead of EntryBooking.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FacilityRegister")

$ws.Range("E2").Value = "KALEB"
$ws.Range("F2").Value = "CHAIRUL"

$ws.Activate() | Out-Null
$ws.Range("F3").Select() | Out-Null
